$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Nearly all code is executed on the OpenCV Frame thread." (Level2)
#    -> "Main Thread:  - That which constructs the app?"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Nearly all code is executed on the OpenCV Frame thread.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Main Thread:  - That which constructs the app?", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Tracking and Feedback" (Heading2 + lastRenderedPageBreak)
#    -> text becomes "UI Thread: - Menu events are handled here", style Level2
#    (Find/Replace on the whole run text also drops the lastRenderedPageBreak
#     marker run since the run gets rewritten.)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Tracking and Feedback",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "UI Thread: - Menu events are handled here", 2) | Out-Null

# Locate that paragraph again (still the 25th paragraph) and fix its style.
$pUi = $d.Paragraphs.Item(25)
$pUi.Style = "Level 2"

# ---------------------------------------------------------------------------
# 3) Insert three new paragraphs after the "UI Thread" paragraph and before
#    "For simplicity, ...":
#      - "OpenCV Frame Thread:  - Camera Frame image is provided and processed."  (Level2)
#      - "OpenGL Thread: - Separate thread on which OpenGL renderers are executed." (Level2, with _GoBack bookmark)
#      - an empty Level2 paragraph
#    Then re-insert a "Tracking and Feedback" Heading2 paragraph right after
#    that, ahead of "For simplicity...".
# ---------------------------------------------------------------------------
$pUi.Range.InsertParagraphAfter() | Out-Null
$pOpenCv = $d.Paragraphs.Item(26)
$pOpenCv.Range.Text = "OpenCV Frame Thread:  - Camera Frame image is provided and processed."

$pOpenCv.Range.InsertParagraphAfter() | Out-Null
$pOpenGl = $d.Paragraphs.Item(27)
$pOpenGl.Range.Text = "OpenGL Thread: - Separate thread on which OpenGL renderers are executed."

$pOpenGl.Range.InsertParagraphAfter() | Out-Null
$pEmpty = $d.Paragraphs.Item(28)
$pEmpty.Style = "Level 2"

$pEmpty.Range.InsertParagraphAfter() | Out-Null
$pTracking = $d.Paragraphs.Item(29)
$pTracking.Range.Text = "Tracking and Feedback"
$pTracking.Style = "Heading 2"

# ---------------------------------------------------------------------------
# 4) Put the "_GoBack" bookmark (collapsed) inside the OpenGL Thread paragraph,
#    right between "...are" and " executed." -- remove it from its old spot
#    (immediately before the Class Diagram drawing) first.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$openGlText = $pOpenGl.Range.Text
$marker = "are"
$idx = $openGlText.IndexOf($marker) + $marker.Length
$bmPos = $pOpenGl.Range.Start + $idx
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
